$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 85.36364133333332
$ws.Range("H2").Value = 256.090924
$ws.Range("I2").Value = 0.832590152283795
$ws.Range("J2").Value = 0.8325901522837948
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 2.300909333333333
$ws.Range("N2").Value = 6.902728
$ws.Range("O2").Value = 0.03776979643482627
$ws.Range("P2").Value = 0.03776979643482627
$ws.Range("Q2").Value = 196.4139990711857
$ws.Range("R2").Value = 1767.725991640672
$ws.Range("S2").Value = 0.03144676056539995
$ws.Range("T2").Value = 0.03144676056539993

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 85.36364133333332
$ws.Range("H3").Value = 256.090924
$ws.Range("I3").Value = 0.832590152283795
$ws.Range("J3").Value = 0.8325901522837948
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 16.96312166666667
$ws.Range("N3").Value = 50.889365
$ws.Range("O3").Value = 0.2784523679257784
$ws.Range("P3").Value = 0.2784523679257784
$ws.Range("Q3").Value = 1448.033833847029
$ws.Range("R3").Value = 13032.30450462326
$ws.Range("S3").Value = 0.2318366994151072
$ws.Range("T3").Value = 0.2318366994151072

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 85.36364133333332
$ws.Range("H4").Value = 256.090924
$ws.Range("I4").Value = 0.832590152283795
$ws.Range("J4").Value = 0.8325901522837948
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 41.655263
$ws.Range("N4").Value = 124.965789
$ws.Range("O4").Value = 0.6837778356393953
$ws.Range("P4").Value = 0.6837778356393953
$ws.Range("Q4").Value = 3555.84493037767
$ws.Range("R4").Value = 32002.60437339903
$ws.Range("S4").Value = 0.5693066923032879
$ws.Range("T4").Value = 0.5693066923032878

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 11.230072
$ws.Range("H5").Value = 33.690216
$ws.Range("I5").Value = 0.1095319647872954
$ws.Range("J5").Value = 0.1095319647872954
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 2.300909333333333
$ws.Range("N5").Value = 6.902728
$ws.Range("O5").Value = 0.03776979643482627
$ws.Range("P5").Value = 0.03776979643482627
$ws.Range("Q5").Value = 25.83937747880533
$ws.Range("R5").Value = 232.554397309248
$ws.Range("S5").Value = 0.004137000013122708
$ws.Range("T5").Value = 0.004137000013122706

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 11.230072
$ws.Range("H6").Value = 33.690216
$ws.Range("I6").Value = 0.1095319647872954
$ws.Range("J6").Value = 0.1095319647872954
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 16.96312166666667
$ws.Range("N6").Value = 50.889365
$ws.Range("O6").Value = 0.2784523679257784
$ws.Range("P6").Value = 0.2784523679257784
$ws.Range("Q6").Value = 190.4970776614267
$ws.Range("R6").Value = 1714.47369895284
$ws.Range("S6").Value = 0.03049943495858539
$ws.Range("T6").Value = 0.03049943495858539

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 11.230072
$ws.Range("H7").Value = 33.690216
$ws.Range("I7").Value = 0.1095319647872954
$ws.Range("J7").Value = 0.1095319647872954
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 41.655263
$ws.Range("N7").Value = 124.965789
$ws.Range("O7").Value = 0.6837778356393953
$ws.Range("P7").Value = 0.6837778356393953
$ws.Range("Q7").Value = 467.791602668936
$ws.Range("R7").Value = 4210.124424020424
$ws.Range("S7").Value = 0.07489552981558732
$ws.Range("T7").Value = 0.07489552981558731

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 5.934092333333333
$ws.Range("H8").Value = 17.802277
$ws.Range("I8").Value = 0.05787788292890966
$ws.Range("J8").Value = 0.05787788292890966
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 2.300909333333333
$ws.Range("N8").Value = 6.902728
$ws.Range("O8").Value = 0.03776979643482627
$ws.Range("P8").Value = 0.03776979643482627
$ws.Range("Q8").Value = 13.65380843462844
$ws.Range("R8").Value = 122.884275911656
$ws.Range("S8").Value = 0.002186035856303625
$ws.Range("T8").Value = 0.002186035856303624

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 5.934092333333333
$ws.Range("H9").Value = 17.802277
$ws.Range("I9").Value = 0.05787788292890966
$ws.Range("J9").Value = 0.05787788292890966
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 16.96312166666667
$ws.Range("N9").Value = 50.889365
$ws.Range("O9").Value = 0.2784523679257784
$ws.Range("P9").Value = 0.2784523679257784
$ws.Range("Q9").Value = 100.6607302315672
$ws.Range("R9").Value = 905.946572084105
$ws.Range("S9").Value = 0.01611623355208588
$ws.Range("T9").Value = 0.01611623355208588

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 5.934092333333333
$ws.Range("H10").Value = 17.802277
$ws.Range("I10").Value = 0.05787788292890966
$ws.Range("J10").Value = 0.05787788292890966
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 41.655263
$ws.Range("N10").Value = 124.965789
$ws.Range("O10").Value = 0.6837778356393953
$ws.Range("P10").Value = 0.6837778356393953
$ws.Range("Q10").Value = 247.1861768112836
$ws.Range("R10").Value = 2224.675591301553
$ws.Range("S10").Value = 0.03957561352052016
$ws.Range("T10").Value = 0.03957561352052015

